$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.114099999999997
$ws.Range("D4").Value = -6.933399999999998
$ws.Range("A11").Value = -21.7836
$ws.Range("A12").Value = -21.5154
$ws.Range("D14").Value = -7.572800000000004
$ws.Range("A15").Value = -21.883
$ws.Range("D26").Value = -8.529900000000003
$ws.Range("A27").Value = -21.82379999999999
$ws.Range("A28").Value = -21.855
$ws.Range("A31").Value = -21.8265
$ws.Range("D31").Value = -8.318100000000001
$ws.Range("A32").Value = -21.3546
$ws.Range("D35").Value = -8.541099999999993
$ws.Range("A36").Value = -21.0779
$ws.Range("D37").Value = -7.868299999999996
$ws.Range("A38").Value = -20.09719999999999
$ws.Range("D39").Value = -7.9907
$ws.Range("D40").Value = -7.911299999999995
$ws.Range("D45").Value = -7.661899999999999
$ws.Range("A46").Value = -21.9054
$ws.Range("D52").Value = -7.527599999999996
$ws.Range("A54").Value = -21.60689999999999
$ws.Range("A55").Value = -22.26840000000001
$ws.Range("A56").Value = -22.12380000000001
$ws.Range("D57").Value = -8.321300000000003
$ws.Range("A67").Value = -21.52519999999998
$ws.Range("A69").Value = -21.67889999999997
$ws.Range("A72").Value = -21.79569999999999
$ws.Range("A73").Value = -19.84329999999999
$ws.Range("D81").Value = -7.120899999999997
$ws.Range("A83").Value = -21.79989999999999
$ws.Range("D83").Value = -8.554899999999996
$ws.Range("A86").Value = -21.95220000000001
$ws.Range("A91").Value = -21.4539
$ws.Range("A93").Value = -21.17739999999999
$ws.Range("A99").Value = -20.07599999999999
$ws.Range("D100").Value = -8.0868
$ws.Range("D102").Value = -8.111399999999996
